$d = $word.ActiveDocument

# Locate "Adobe Photoshop" (bold run, part of the list of design tools) and
# collapse the found range to its end so we can insert right after it.
$rng = $d.Content
$found = $rng.Find.Execute("Adobe Photoshop")
if (-not $found) {
    throw "Could not find 'Adobe Photoshop' in the document"
}
$rng.Collapse(0)

# Insert ", Adobe Illustrator" right after "Adobe Photoshop" - the
# insertion point sits inside the existing bold run, so the new text
# naturally inherits the bold formatting. This turns
#   "Figma, Adobe Photoshop and GIMP for 'Design' phase,"
# into
#   "Figma, Adobe Photoshop, Adobe Illustrator and GIMP for 'Design' phase,"
$rng.InsertAfter(", Adobe Illustrator")

# Word keeps the hidden "_GoBack" bookmark at the location of the most
# recent edit, so it now sits between the inserted ", " and "Adobe
# Illustrator" rather than at the old, now-empty paragraph near the end of
# the section. Re-create it there (adding it elsewhere automatically moves
# it from its previous location).
$full = $d.Content.Text
$idx = $full.IndexOf("Adobe Photoshop, ") + ("Adobe Photoshop, ").Length
$bmRange = $d.Range($idx, $idx)
$d.Bookmarks.Add("_GoBack", $bmRange)
